# Adds a new "First Time Reading?" boolean column (I) to the Completed sheet,
# records whether each book had been read before, and refreshes the sheet
# view (zoom / scroll / selection) to match how the author last left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Completed")

# --- New column header -------------------------------------------------
$ws.Range("I1").Value = "First Time Reading?"

# --- Per-row "read it before?" flags (row 2 .. row 104) ---------------
$values = @($true,$true,$true,$true,$true,$false,$true,$true,$true,$false,$true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$false,$true,$true,$true,$false,$true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$false,$true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$true,$false,$false,$true,$true,$true,$true,$true,$true,$false,$true,$false,$true,$true,$false,$false,$false,$true)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i]
}

# --- Column width for the new column -----------------------------------
$ws.Columns.Item(9).ColumnWidth = 17.08984375

# --- Restore the view state (zoom + selection) the author left it in ---
$win = $excel.ActiveWindow
$win.Zoom = 90
$ws.Range("I105").Select()

"done"
